# Fill in the previously-empty "Stronger penalty + BO5" rows (19 and 20) of
# the seed-investigation report with their measured accuracy numbers, the
# Diff. formula, and the "BO5" label - matching the style already used by
# the rest of the table (thin border, centered alignment; row 20's
# "Incremental" column keeps the 3-decimal numeric display).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 19 ---------------------------------------------------------------
$ws.Range("E19").Value = 0.772
$ws.Range("F19").Value = 0.773
$ws.Range("G19").Formula = "=E19-F19"
$ws.Range("H19").Value = "BO5"

# --- Row 20 -----------------------------------------------------------------
$ws.Range("E20").Value = 0.736
$ws.Range("F20").Value = 0.78
$ws.Range("G20").Formula = "=E20-F20"
$ws.Range("H20").Value = "BO5"

# --- Formatting: thin border + centered alignment, applied a whole row at a
# time (batches the underlying style writes instead of touching cell by
# cell) ----------------------------------------------------------------------
$ws.Range("E19:H19").HorizontalAlignment = -4108
$ws.Range("E19:H19").VerticalAlignment = -4108
$ws.Range("E20:H20").HorizontalAlignment = -4108
$ws.Range("E20:H20").VerticalAlignment = -4108

# F20 keeps the 3-decimal numeric display used elsewhere for "Incremental"
# accuracy values
$ws.Range("F20").NumberFormat = "0.000"

# --- Selection / view --------------------------------------------------------
$ws.Range("G20").Select()

# --- Page setup -------------------------------------------------------------
$ws.PageSetup.Orientation = 1
